# NYPD CompStat 106th Precinct weekly report refresh:
#  - bump the report "Volume/Number" and the covered week date range
#  - refresh all weekly/28-day/YTD crime counts and their % change columns
#    for the newly collected week of crime data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -----------------------------------------------
# "Volume 30   Number  20" -> "Volume 30   Number  22"
$ws.Range("A8").Value2 = "Volume 30   Number  22"
# "Report Covering the Week  5/15/2023  Through  5/21/2023"
#   -> "Report Covering the Week  5/29/2023  Through  6/4/2023"
$ws.Range("C9").Value2 = "Report Covering the Week  5/29/2023  Through  6/4/2023"

# --- Fix up cells whose content type changes (number <-> dash/n-a text) -
# Some cells flip between a plain numeric value and the report's special
# placeholder text ("0" / "***.*"). Grab number-format/style from a
# same-shaped neighboring cell before writing the new content so the
# underlying style index matches what the rest of the column uses.

# G14, F30 become the literal "0" placeholder (style like C14)
$ws.Range("C14").Copy($ws.Range("G14")) | Out-Null
$ws.Range("C14").Copy($ws.Range("F30")) | Out-Null
# H14 becomes the "***.*" placeholder (style like E14)
$ws.Range("E14").Copy($ws.Range("H14")) | Out-Null

# C27, D28, D29 become plain numbers (style like I26)
$ws.Range("I26").Copy($ws.Range("C27")) | Out-Null
$ws.Range("I26").Copy($ws.Range("D28")) | Out-Null
$ws.Range("I26").Copy($ws.Range("D29")) | Out-Null
# E28, E29 become plain % numbers (style like M28/M29)
$ws.Range("M28").Copy($ws.Range("E28")) | Out-Null
$ws.Range("M29").Copy($ws.Range("E29")) | Out-Null

# Now assign the actual new values for those type-changed cells
$ws.Range("G14").Value2 = "0"
$ws.Range("H14").Value2 = "***.*"
$ws.Range("F30").Value2 = "0"
$ws.Range("C27").Value2 = 1
$ws.Range("D28").Value2 = 1
$ws.Range("D29").Value2 = 1
$ws.Range("E28").Value2 = -100
$ws.Range("E29").Value2 = -100

# --- Refresh all the other weekly/28-day/YTD figures and % changes -----
$ws.Range("N14").Value2 = -92.857142857142
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 4
$ws.Range("H15").Value2 = -75
$ws.Range("J15").Value2 = 11
$ws.Range("K15").Value2 = -18.181818181818
$ws.Range("N15").Value2 = -40
$ws.Range("C16").Value2 = 3
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = 200
$ws.Range("F16").Value2 = 10
$ws.Range("G16").Value2 = 11
$ws.Range("H16").Value2 = -9.090909090909
$ws.Range("I16").Value2 = 82
$ws.Range("J16").Value2 = 81
$ws.Range("K16").Value2 = 1.234567901234
$ws.Range("L16").Value2 = 12.328767123287
$ws.Range("M16").Value2 = -14.583333333333
$ws.Range("N16").Value2 = -77.717391304347
$ws.Range("C17").Value2 = 11
$ws.Range("E17").Value2 = 83.333333333333
$ws.Range("F17").Value2 = 30
$ws.Range("G17").Value2 = 25
$ws.Range("H17").Value2 = 20
$ws.Range("I17").Value2 = 137
$ws.Range("J17").Value2 = 119
$ws.Range("K17").Value2 = 15.126050420168
$ws.Range("L17").Value2 = 25.688073394495
$ws.Range("M17").Value2 = 124.590163934426
$ws.Range("N17").Value2 = -11.612903225806
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 6
$ws.Range("F18").Value2 = 6
$ws.Range("G18").Value2 = 14
$ws.Range("H18").Value2 = -57.142857142857
$ws.Range("I18").Value2 = 42
$ws.Range("J18").Value2 = 59
$ws.Range("K18").Value2 = -28.813559322033
$ws.Range("L18").Value2 = 5
$ws.Range("M18").Value2 = -64.406779661017
$ws.Range("N18").Value2 = -91.446028513238
$ws.Range("C19").Value2 = 15
$ws.Range("D19").Value2 = 15
$ws.Range("F19").Value2 = 54
$ws.Range("G19").Value2 = 49
$ws.Range("H19").Value2 = 10.204081632653
$ws.Range("I19").Value2 = 261
$ws.Range("J19").Value2 = 302
$ws.Range("K19").Value2 = -13.576158940397
$ws.Range("L19").Value2 = 36.649214659685
$ws.Range("M19").Value2 = 86.428571428571
$ws.Range("N19").Value2 = 9.663865546218
$ws.Range("C20").Value2 = 8
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = 300
$ws.Range("F20").Value2 = 14
$ws.Range("G20").Value2 = 9
$ws.Range("H20").Value2 = 55.555555555555
$ws.Range("I20").Value2 = 90
$ws.Range("J20").Value2 = 107
$ws.Range("K20").Value2 = -15.887850467289
$ws.Range("L20").Value2 = 55.172413793103
$ws.Range("M20").Value2 = -26.829268292682
$ws.Range("N20").Value2 = -93.715083798882
$ws.Range("C21").Value2 = 40
$ws.Range("D21").Value2 = 31
$ws.Range("E21").Value2 = 29.032258064516
$ws.Range("F21").Value2 = 116
$ws.Range("G21").Value2 = 112
$ws.Range("H21").Value2 = 3.571428571428
$ws.Range("I21").Value2 = 622
$ws.Range("J21").Value2 = 680
$ws.Range("K21").Value2 = -8.529411764705
$ws.Range("L21").Value2 = 29.853862212943
$ws.Range("M21").Value2 = 13.090909090909
$ws.Range("N21").Value2 = -77.073350534463
$ws.Range("F22").Value2 = 1
$ws.Range("L22").Value2 = 14.285714285714
$ws.Range("M22").Value2 = -27.272727272727
$ws.Range("C24").Value2 = 20
$ws.Range("D24").Value2 = 17
$ws.Range("E24").Value2 = 17.647058823529
$ws.Range("F24").Value2 = 113
$ws.Range("G24").Value2 = 111
$ws.Range("H24").Value2 = 1.801801801801
$ws.Range("I24").Value2 = 542
$ws.Range("J24").Value2 = 594
$ws.Range("K24").Value2 = -8.754208754208
$ws.Range("L24").Value2 = 59.411764705882
$ws.Range("M24").Value2 = 100.740740740741
$ws.Range("C25").Value2 = 15
$ws.Range("D25").Value2 = 12
$ws.Range("E25").Value2 = 25
$ws.Range("F25").Value2 = 47
$ws.Range("G25").Value2 = 51
$ws.Range("H25").Value2 = -7.843137254901
$ws.Range("I25").Value2 = 234
$ws.Range("J25").Value2 = 193
$ws.Range("K25").Value2 = 21.243523316062
$ws.Range("L25").Value2 = 32.203389830508
$ws.Range("M25").Value2 = 12.5
$ws.Range("D26").Value2 = 2
$ws.Range("F26").Value2 = 1
$ws.Range("G26").Value2 = 6
$ws.Range("H26").Value2 = -83.333333333333
$ws.Range("J26").Value2 = 19
$ws.Range("K26").Value2 = -31.578947368421
$ws.Range("D27").Value2 = 1
$ws.Range("E27").Value2 = 0
$ws.Range("F27").Value2 = 3
$ws.Range("G27").Value2 = 7
$ws.Range("H27").Value2 = -57.142857142857
$ws.Range("I27").Value2 = 24
$ws.Range("J27").Value2 = 26
$ws.Range("K27").Value2 = -7.692307692307
$ws.Range("L27").Value2 = 26.315789473684
$ws.Range("G28").Value2 = 2
$ws.Range("J28").Value2 = 6
$ws.Range("K28").Value2 = -33.333333333333
$ws.Range("L28").Value2 = -20
$ws.Range("G29").Value2 = 2
$ws.Range("J29").Value2 = 6
$ws.Range("K29").Value2 = -66.666666666666
$ws.Range("L29").Value2 = -50
